$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Slamf7"
$ws.Range("C2").Value2 = "Slamf7"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.007936333333333333
$ws.Range("H2").Value2 = 0.023809
$ws.Range("I2").Value2 = 0.768081811729789
$ws.Range("J2").Value2 = 0.768081811729789
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 0.007936333333333333
$ws.Range("N2").Value2 = 0.023809
$ws.Range("O2").Value2 = 0.768081811729789
$ws.Range("P2").Value2 = 0.768081811729789
$ws.Range("Q2").Value2 = 0.00006298538677777778
$ws.Range("R2").Value2 = 0.000566868481
$ws.Range("S2").Value2 = 0.589949669510115
$ws.Range("T2").Value2 = 0.589949669510115

# Row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Slamf7"
$ws.Range("C3").Value2 = "Slamf7"
$ws.Range("D3").Value2 = "MuSCs"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.007936333333333333
$ws.Range("H3").Value2 = 0.023809
$ws.Range("I3").Value2 = 0.768081811729789
$ws.Range("J3").Value2 = 0.768081811729789
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0.3333333333333333
$ws.Range("M3").Value2 = 0.002396333333333333
$ws.Range("N3").Value2 = 0.007189
$ws.Range("O3").Value2 = 0.231918188270211
$ws.Range("P3").Value2 = 0.231918188270211
$ws.Range("Q3").Value2 = 0.00001901810011111111
$ws.Range("R3").Value2 = 0.000171162901
$ws.Range("S3").Value2 = 0.1781321422196739
$ws.Range("T3").Value2 = 0.1781321422196739

# Row 4
$ws.Range("A4").Value2 = "MuSCs"
$ws.Range("B4").Value2 = "Slamf7"
$ws.Range("C4").Value2 = "Slamf7"
$ws.Range("D4").Value2 = "ECs"
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.002396333333333333
$ws.Range("H4").Value2 = 0.007189
$ws.Range("I4").Value2 = 0.231918188270211
$ws.Range("J4").Value2 = 0.231918188270211
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.007936333333333333
$ws.Range("N4").Value2 = 0.023809
$ws.Range("O4").Value2 = 0.768081811729789
$ws.Range("P4").Value2 = 0.768081811729789
$ws.Range("Q4").Value2 = 0.00001901810011111111
$ws.Range("R4").Value2 = 0.000171162901
$ws.Range("S4").Value2 = 0.1781321422196739
$ws.Range("T4").Value2 = 0.1781321422196739

# Row 5
$ws.Range("A5").Value2 = "MuSCs"
$ws.Range("B5").Value2 = "Slamf7"
$ws.Range("C5").Value2 = "Slamf7"
$ws.Range("D5").Value2 = "MuSCs"
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.002396333333333333
$ws.Range("H5").Value2 = 0.007189
$ws.Range("I5").Value2 = 0.231918188270211
$ws.Range("J5").Value2 = 0.231918188270211
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.002396333333333333
$ws.Range("N5").Value2 = 0.007189
$ws.Range("O5").Value2 = 0.231918188270211
$ws.Range("P5").Value2 = 0.231918188270211
$ws.Range("Q5").Value2 = 0.000005742413444444444
$ws.Range("R5").Value2 = 0.000051681721
$ws.Range("S5").Value2 = 0.05378604605053702
$ws.Range("T5").Value2 = 0.05378604605053702

Write-Output "Edit complete"